# Apply the "Update with Correct Forecast output" edit:
#  - Forecast Comparison sheet: insert a new "Week_Start_Date" column (B),
#    shorten the Week labels (W01 -> W1 ...), correct two MyForecast values,
#    and store is_holiday_week as boolean values.
#  - Summary sheet: refresh the two summary metrics that depend on the
#    corrected forecast (Total Forecast (4 Weeks) and Max Forecast).

$wb2 = $excel.ActiveWorkbook
$ws  = $wb2.Worksheets.Item("Forecast Comparison")
$sum = $wb2.Worksheets.Item("Summary")

# --- Insert the new "Week_Start_Date" column between Week (A) and ASIN (B) ---
$ws.Columns.Item(2).Insert()
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# --- Week_Start_Date values (stored as text, matching the week-start dates) ---
$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)
for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    # Leading apostrophe forces the date-looking string to stay text.
    $ws.Cells.Item($row, 2).Value = "'" + $weekStartDates[$i]
}

# --- Shorten the Week labels: "W01".."W16" -> "W1".."W16" ---
for ($i = 1; $i -le 16; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "W" + $i
}

# --- Correct the MyForecast values (now column D after the insert) ---
$ws.Cells.Item(2, 4).Value = 171
$ws.Cells.Item(9, 4).Value = 183

# --- is_holiday_week (now column J) should be stored as boolean values ---
$ws.Range("J2:J17").Value = $false

# --- Update the dependent Summary metrics (kept as text, like the rest of
#     the Summary sheet's Value column) ---
$sum.Cells.Item(11, 2).Value = "'597"
$sum.Cells.Item(12, 2).Value = "'183"
